# TOD-E norms run, POM rescale, 24 cell demo strat
#
# Splits the old combined "7.0-9.3" raw->standard-score lookup tab into four
# narrower age-band tabs (7.0-7.5 / 7.6-7.11 / 8.0-8.5 / 8.6-9.3), rescaling
# the "ss" column for each band.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the existing "7.0-9.3" sheet to "7.0-7.5" and rescale its ss
#    column (raw 1..32 in col A is unchanged; only col B "ss" values move).
# ---------------------------------------------------------------------------
$sheet6 = $wb.Worksheets.Item(6)
$sheet6.Name = "7.0-7.5"

$sheet6_ss = @(54,57,60,62,65,67,70,72,75,77,79,82,84,87,89,91,94,96,99,102,104,107,110,112,115,118,121,125,128,130,130,130)
for ($i = 0; $i -lt $sheet6_ss.Length; $i++) {
    $sheet6.Cells.Item($i + 2, 2).Value = $sheet6_ss[$i]
}

# ---------------------------------------------------------------------------
# Helper data for the three brand-new tabs. Each tab has the same shape as
# the others: header row ("raw","ss") then raw=1..32 with a rescaled ss.
# ---------------------------------------------------------------------------
$sheet7_ss = @(52,54,57,59,61,64,66,68,71,73,75,78,80,82,84,87,89,92,94,96,99,101,104,107,110,112,115,119,122,126,130,130)
$sheet8_ss = @(49,52,54,57,59,61,63,66,68,70,72,75,77,79,81,84,86,88,91,93,95,98,101,103,106,109,112,115,119,123,129,130)
$sheet9_ss = @(48,50,53,55,57,59,62,64,66,68,71,73,75,77,80,82,84,87,89,92,94,97,100,103,107,110,115,120,130,130,130,130)

function Add-NormsTab($Name, $AfterSheet, $SsValues) {
    $ws = $wb.Worksheets.Add($null, $AfterSheet)
    $ws.Name = $Name

    $ws.Cells.Item(1, 1).Value = "raw"
    $ws.Cells.Item(1, 2).Value = "ss"
    $header = $ws.Range("A1:B1")
    $header.Font.Bold = $true
    $header.HorizontalAlignment = -4108

    for ($i = 0; $i -lt $SsValues.Length; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i + 1
        $ws.Cells.Item($i + 2, 2).Value = $SsValues[$i]
    }

    return $ws
}

# ---------------------------------------------------------------------------
# 2) Insert the three new tabs right after "7.0-7.5", in order, each one
#    after the one before it so the final left-to-right order is:
#    ... 7.0-7.5, 7.6-7.11, 8.0-8.5, 8.6-9.3
# ---------------------------------------------------------------------------
$sheet7 = Add-NormsTab "7.6-7.11" $sheet6 $sheet7_ss
$sheet8 = Add-NormsTab "8.0-8.5" $sheet7 $sheet8_ss
$sheet9 = Add-NormsTab "8.6-9.3" $sheet8 $sheet9_ss

# Keep the originally-selected first tab active.
$wb.Worksheets.Item(1).Activate()
